$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 69.44
$ws.Range("I15").Value = 69.44
$ws.Range("K15").Value = 208.32
$ws.Range("M15").Value = -39.31999999999999
$ws.Range("H138").Value = 2026.6061
$ws.Range("I138").Value = 1249.1794
$ws.Range("J138").Value = 3149.5557
$ws.Range("K138").Value = 3747.5382
$ws.Range("L138").Value = 9448.667099999999
$ws.Range("M138").Value = 1392.4618
$ws.Range("N138").Value = -19728.6671
$ws.Range("H139").Value = 64750
$ws.Range("J139").Value = 64750
$ws.Range("L139").Value = 64750
$ws.Range("N139").Value = -75030

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 813.4706
$ws.Range("I2").Value = 739.375
$ws.Range("J2").Value = 1999
$ws.Range("K2").Value = 739.375
$ws.Range("L2").Value = 1999
$ws.Range("M2").Value = -626.375
$ws.Range("N2").Value = -2225
$ws.Range("H45").Value = 2228.8462
$ws.Range("I45").Value = 1089.4
$ws.Range("J45").Value = 2941
$ws.Range("K45").Value = 1089.4
$ws.Range("L45").Value = 2941
$ws.Range("M45").Value = -712.4000000000001
$ws.Range("N45").Value = -3695
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("N109").Value = 0
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("N111").Value = 0
$ws.Range("H113").Value = 41993.332
$ws.Range("J113").Value = 41993.332
$ws.Range("L113").Value = 41993.332
$ws.Range("N113").Value = -50671.332
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("N114").Value = 0
$ws.Range("H116").Value = 813.4706
$ws.Range("I116").Value = 739.375
$ws.Range("J116").Value = 1999
$ws.Range("K116").Value = 739.375
$ws.Range("L116").Value = 1999
$ws.Range("M116").Value = 1554.625
$ws.Range("N116").Value = -6587
$ws.Range("H117").Value = 43000
$ws.Range("J117").Value = 43000
$ws.Range("L117").Value = 43000
$ws.Range("N117").Value = -52178
$ws.Range("H119").Value = 48000
$ws.Range("J119").Value = 48000
$ws.Range("L119").Value = 48000
$ws.Range("N119").Value = -57676
$ws.Range("H122").Value = 859.73334
$ws.Range("I122").Value = 806.8570999999999
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 2420.5713
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = 29.42870000000039
$ws.Range("N122").Value = -9700
$ws.Range("H124").Value = 34000
$ws.Range("J124").Value = 34000
$ws.Range("L124").Value = 34000
$ws.Range("N124").Value = -43820
$ws.Range("H125").Value = 48000
$ws.Range("J125").Value = 48000
$ws.Range("L125").Value = 48000
$ws.Range("N125").Value = -57840
$ws.Range("H128").Value = 48000
$ws.Range("J128").Value = 48000
$ws.Range("L128").Value = 48000
$ws.Range("N128").Value = -57960
$ws.Range("H129").Value = 10709
$ws.Range("I129").Value = 10709
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 10709
$ws.Range("M129").Value = -5709
$ws.Range("N129").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 813.4706
$ws.Range("I3").Value = 739.375
$ws.Range("J3").Value = 1999
$ws.Range("K3").Value = 739.375
$ws.Range("L3").Value = 1999
$ws.Range("M3").Value = -625.375
$ws.Range("N3").Value = -2227
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("N109").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3877.3928
$ws.Range("I86").Value = 3371.2942
$ws.Range("J86").Value = 4659.5454
$ws.Range("K86").Value = 3371.2942
$ws.Range("L86").Value = 4659.5454
$ws.Range("M86").Value = -2248.2942
$ws.Range("N86").Value = -6905.5454
$ws.Range("H89").Value = 3877.3928
$ws.Range("I89").Value = 3371.2942
$ws.Range("J89").Value = 4659.5454
$ws.Range("K89").Value = 16856.471
$ws.Range("L89").Value = 23297.727
$ws.Range("M89").Value = -11240.471
$ws.Range("N89").Value = -34529.727
$ws.Range("H132").Value = 4964.9546
$ws.Range("I132").Value = 5609.8335
$ws.Range("K132").Value = 16829.5005
$ws.Range("M132").Value = -14299.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 984.8
$ws.Range("J127").Value = 984.8
$ws.Range("L127").Value = 2954.4
$ws.Range("N127").Value = -12874.4
$ws.Range("H130").Value = 1333.3334
$ws.Range("H131").Value = 868.4651
$ws.Range("J131").Value = 1002.4286
$ws.Range("L131").Value = 3007.2858
$ws.Range("N131").Value = -13087.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 44999
$ws.Range("J42").Value = 44999
$ws.Range("L42").Value = 44999
$ws.Range("N42").Value = -45969
$ws.Range("H102").Value = 6633.1665
$ws.Range("I102").Value = 8668
$ws.Range("J102").Value = 3784.4
$ws.Range("K102").Value = 8668
$ws.Range("L102").Value = 3784.4
$ws.Range("M102").Value = -7046
$ws.Range("N102").Value = -7028.4
$ws.Range("H115").Value = 44999
$ws.Range("J115").Value = 44999
$ws.Range("L115").Value = 44999
$ws.Range("N115").Value = -47349
$ws.Range("H119").Value = 47600
$ws.Range("J119").Value = 47600
$ws.Range("L119").Value = 47600
$ws.Range("N119").Value = -57276
$ws.Range("H122").Value = 1794.3334
$ws.Range("I122").Value = 1208.2142
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 3624.6426
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -1174.6426
$ws.Range("N122").Value = -34900
$ws.Range("H123").Value = 40326
$ws.Range("J123").Value = 40326
$ws.Range("L123").Value = 40326
$ws.Range("N123").Value = -45226
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("N124").Value = 0
$ws.Range("H127").Value = 31442
$ws.Range("J127").Value = 31442
$ws.Range("L127").Value = 31442
$ws.Range("N127").Value = -41362
$ws.Range("H136").Value = 12803.473
$ws.Range("J136").Value = 12803.473
$ws.Range("L136").Value = 38410.419
$ws.Range("N136").Value = -43510.419

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1661.6842
$ws.Range("I7").Value = 1604.6
$ws.Range("J7").Value = 1875.75
$ws.Range("K7").Value = 1604.6
$ws.Range("L7").Value = 1875.75
$ws.Range("M7").Value = -1492.6
$ws.Range("N7").Value = -2099.75
$ws.Range("H40").Value = 1957.9286
$ws.Range("I40").Value = 1344
$ws.Range("J40").Value = 2571.8572
$ws.Range("K40").Value = 1344
$ws.Range("L40").Value = 2571.8572
$ws.Range("M40").Value = -1208
$ws.Range("N40").Value = -2843.8572
$ws.Range("H122").Value = 61929.53
$ws.Range("J122").Value = 3781.8
$ws.Range("L122").Value = 11345.4
$ws.Range("N122").Value = -16245.4
$ws.Range("H123").Value = 24571.428
$ws.Range("J123").Value = 24571.428
$ws.Range("L123").Value = 24571.428
$ws.Range("N123").Value = -34371.428
$ws.Range("H126").Value = 1661.6842
$ws.Range("I126").Value = 1604.6
$ws.Range("J126").Value = 1875.75
$ws.Range("K126").Value = 4813.799999999999
$ws.Range("L126").Value = 5627.25
$ws.Range("M126").Value = -2343.799999999999
$ws.Range("N126").Value = -10567.25
$ws.Range("H132").Value = 19882.416
$ws.Range("I132").Value = 13901.6
$ws.Range("J132").Value = 24154.428
$ws.Range("K132").Value = 41704.8
$ws.Range("L132").Value = 72463.284
$ws.Range("M132").Value = -39174.8
$ws.Range("N132").Value = -77523.284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1002.619
$ws.Range("I122").Value = 958.8889
$ws.Range("J122").Value = 1265
$ws.Range("K122").Value = 2876.6667
$ws.Range("L122").Value = 3795
$ws.Range("M122").Value = -426.6667000000002
$ws.Range("N122").Value = -8695
$ws.Range("H132").Value = 2567.831
$ws.Range("I132").Value = 2952.889
$ws.Range("K132").Value = 8858.667000000001
$ws.Range("M132").Value = -6328.667000000001
